$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B24").Value = 21
$ws.Range("C24").Value = "Comfort edit dates"
$ws.Range("C25").Select() | Out-Null
